$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the daily dataset. It lands at row 286,
# pushing the existing rows 286-313 down to 287-314 (dimension grows to R314).
$ws.Rows.Item(286).Insert()

# Populate the newly inserted row 286 with the new observation.
$ws.Range("A286").Value = 7
$ws.Range("B286").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C286").Value = "Ñuble"
$ws.Range("D286").Value = 44783
$ws.Range("E286").Value = 16
$ws.Range("F286").Value = 100114013
$ws.Range("G286").Value = "Zanahoria"
$ws.Range("H286").Value = "Sin especificar"
$ws.Range("I286").Value = "Primera"
$ws.Range("J286").Value = 120
$ws.Range("K286").Value = 8000
$ws.Range("L286").Value = 9000
$ws.Range("M286").Value = 8500
$ws.Range("N286").Value = "$/saco 20 kilos"
$ws.Range("O286").Value = "Región de Ñuble"
$ws.Range("P286").Value = 425
$ws.Range("Q286").Value = 20
$ws.Range("R286").Value = "Hortaliza"
